# This sheet holds one weekly price-report data table (header in row 1,
# data rows starting at row 2). A new weekly observation was inserted
# right before the existing row 14, which pushes every row from 14
# downward down by one position; the previously-last row (130) becomes
# the new row 131.
#
# We therefore walk the rows bottom-up (so we never overwrite a source
# row before it has been copied) shifting row (r-1) into row r for every
# row from the new last row (131) back up to row 15. Row 14 keeps all of
# its original values except for the date (column D), which receives the
# brand-new observation's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstShiftRow = 15
$lastRow = 131
$lastCol = 18   # columns A..R

for ($r = $lastRow; $r -ge $firstShiftRow; $r--) {
    $src = $r - 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $v = $ws.Cells.Item($src, $c).Value2()
        $ws.Cells.Item($r, $c).Value = $v
    }
}

# The brand-new row 131 (column D, the "Fecha" column) needs the same
# date number format that every other row in column D already uses,
# since the cell did not previously exist.
$ws.Cells.Item(131, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 14 only gets a new date value for the freshly inserted observation;
# every other field on that row is unchanged.
$ws.Cells.Item(14, 4).Value = 44649
